$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the existing header formatting (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-9
$saveValues = @(0, 1, 1, 1, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
